$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 30003
$ws.Range("E2").Value = 1276
$ws.Range("F2").Value = 1276
$ws.Range("G2").Value = 1029
$ws.Range("H2").Value = 737
$ws.Range("I2").Value = 736
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 28591
$ws.Range("L2").Value = 18107
$ws.Range("M2").Value = 10484
$ws.Range("N2").Value = 10464
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 1401
$ws.Range("Q2").Value = 353
$ws.Range("R2").Value = -618
$ws.Range("S2").Value = 230
$ws.Range("T2").Value = 780
$ws.Range("U2").Value = -427
$ws.Range("V2").Value = 13514
$ws.Range("W2").Value = 4.25
$ws.Range("X2").Value = 2.46
$ws.Range("Y2").Value = 7.23
$ws.Range("Z2").Value = 2.66
$ws.Range("AA2").Value = 172.72
$ws.Range("AB2").Value = 649.37
$ws.Range("AC2").Value = 2626
$ws.Range("AD2").Value = 9.539999999999999
$ws.Range("AE2").Value = 37355
$ws.Range("AF2").Value = 0.67
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 2.4
$ws.Range("AI2").Value = 22.84
$ws.Range("AJ2").Value = 28024278
$ws.Range("D3").Value = 28197
$ws.Range("E3").Value = 1111
$ws.Range("F3").Value = 1111
$ws.Range("G3").Value = 798
$ws.Range("H3").Value = 517
$ws.Range("I3").Value = 516
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 26076
$ws.Range("L3").Value = 15146
$ws.Range("M3").Value = 10930
$ws.Range("N3").Value = 10910
$ws.Range("O3").Value = 20
$ws.Range("P3").Value = 1401
$ws.Range("Q3").Value = 3396
$ws.Range("R3").Value = -849
$ws.Range("S3").Value = -2456
$ws.Range("T3").Value = 903
$ws.Range("U3").Value = 2494
$ws.Range("V3").Value = 11249
$ws.Range("W3").Value = 3.94
$ws.Range("X3").Value = 1.83
$ws.Range("Y3").Value = 4.83
$ws.Range("Z3").Value = 1.89
$ws.Range("AA3").Value = 138.57
$ws.Range("AB3").Value = 674.42
$ws.Range("AC3").Value = 1841
$ws.Range("AD3").Value = 14.97
$ws.Range("AE3").Value = 38947
$ws.Range("AF3").Value = 0.71
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.81
$ws.Range("AI3").Value = 27.16
$ws.Range("AJ3").Value = 28024278
$ws.Range("D4").Value = 28318
$ws.Range("E4").Value = 2178
$ws.Range("F4").Value = 2178
$ws.Range("G4").Value = 1822
$ws.Range("H4").Value = 1377
$ws.Range("I4").Value = 1376
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 27149
$ws.Range("L4").Value = 14801
$ws.Range("M4").Value = 12348
$ws.Range("N4").Value = 12328
$ws.Range("O4").Value = 21
$ws.Range("P4").Value = 1401
$ws.Range("Q4").Value = 2463
$ws.Range("R4").Value = -1008
$ws.Range("S4").Value = -1287
$ws.Range("T4").Value = 987
$ws.Range("U4").Value = 1475
$ws.Range("V4").Value = 10089
$ws.Range("W4").Value = 7.69
$ws.Range("X4").Value = 4.86
$ws.Range("Y4").Value = 11.84
$ws.Range("Z4").Value = 5.17
$ws.Range("AA4").Value = 119.86
$ws.Range("AB4").Value = 770.52
$ws.Range("AC4").Value = 4910
$ws.Range("AD4").Value = 8.25
$ws.Range("AE4").Value = 44008
$ws.Range("AF4").Value = 0.92
$ws.Range("AG4").Value = 700
$ws.Range("AH4").Value = 1.73
$ws.Range("AI4").Value = 14.25
$ws.Range("AJ4").Value = 28024278
$ws.Range("D5").Value = 29450
$ws.Range("E5").Value = 2411
$ws.Range("F5").Value = 2411
$ws.Range("G5").Value = 2048
$ws.Range("H5").Value = 1507
$ws.Range("I5").Value = 1506
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 27295
$ws.Range("L5").Value = 14204
$ws.Range("M5").Value = 13091
$ws.Range("N5").Value = 13071
$ws.Range("O5").Value = 21
$ws.Range("P5").Value = 1401
$ws.Range("Q5").Value = 1374
$ws.Range("R5").Value = -877
$ws.Range("S5").Value = -565
$ws.Range("T5").Value = 818
$ws.Range("U5").Value = 556
$ws.Range("V5").Value = 9243
$ws.Range("W5").Value = 8.19
$ws.Range("X5").Value = 5.12
$ws.Range("Y5").Value = 11.86
$ws.Range("Z5").Value = 5.54
$ws.Range("AA5").Value = 108.5
$ws.Range("AB5").Value = 840.26
$ws.Range("AC5").Value = 5372
$ws.Range("AD5").Value = 8.84
$ws.Range("AE5").Value = 46660
$ws.Range("AF5").Value = 1.02
$ws.Range("AG5").Value = 800
$ws.Range("AH5").Value = 1.68
$ws.Range("AI5").Value = 14.88
$ws.Range("AJ5").Value = 28024278
$ws.Range("D6").Value = 27745
$ws.Range("E6").Value = 1075
$ws.Range("F6").Value = 1075
$ws.Range("G6").Value = 874
$ws.Range("H6").Value = 621
$ws.Range("I6").Value = 620
$ws.Range("K6").Value = 26681
$ws.Range("L6").Value = 12888
$ws.Range("M6").Value = 13793
$ws.Range("N6").Value = 13772
$ws.Range("P6").Value = 1401
$ws.Range("Q6").Value = 812
$ws.Range("R6").Value = -751
$ws.Range("S6").Value = -85
$ws.Range("T6").Value = 738
$ws.Range("U6").Value = 74
$ws.Range("V6").Value = 9253
$ws.Range("W6").Value = 3.88
$ws.Range("X6").Value = 2.24
$ws.Range("Y6").Value = 4.62
$ws.Range("Z6").Value = 2.3
$ws.Range("AA6").Value = 93.44
$ws.Range("AB6").Value = 882.96
$ws.Range("AC6").Value = 2212
$ws.Range("AD6").Value = 12.36
$ws.Range("AE6").Value = 49164
$ws.Range("AF6").Value = 0.5600000000000001
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 2.19
$ws.Range("AI6").Value = 27.11
$ws.Range("AJ6").Value = 28024278
$ws.Range("D7").Value = 24136
$ws.Range("E7").Value = 397
$ws.Range("G7").Value = 259
$ws.Range("H7").Value = 163
$ws.Range("I7").Value = 173
$ws.Range("K7").Value = 26718
$ws.Range("L7").Value = 12817
$ws.Range("M7").Value = 13902
$ws.Range("N7").Value = 13872
$ws.Range("P7").Value = 1400
$ws.Range("Q7").Value = 1272
$ws.Range("R7").Value = -869
$ws.Range("S7").Value = -388
$ws.Range("T7").Value = 821
$ws.Range("U7").Value = 225
$ws.Range("W7").Value = 1.65
$ws.Range("X7").Value = 0.68
$ws.Range("Y7").Value = 1.26
$ws.Range("Z7").Value = 0.61
$ws.Range("AA7").Value = 92.2
$ws.Range("AC7").Value = 619
$ws.Range("AD7").Value = 33.2
$ws.Range("AE7").Value = 49522
$ws.Range("AF7").Value = 0.41
$ws.Range("AG7").Value = 556
$ws.Range("AH7").Value = 2.7
$ws.Range("AI7").Value = 89.75
$ws.Range("D8").Value = 25593
$ws.Range("E8").Value = 946
$ws.Range("G8").Value = 744
$ws.Range("H8").Value = 547
$ws.Range("I8").Value = 554
$ws.Range("K8").Value = 27331
$ws.Range("L8").Value = 12984
$ws.Range("M8").Value = 14346
$ws.Range("N8").Value = 14314
$ws.Range("P8").Value = 1400
$ws.Range("Q8").Value = 1050
$ws.Range("R8").Value = -763
$ws.Range("S8").Value = -120
$ws.Range("T8").Value = 742
$ws.Range("U8").Value = 142
$ws.Range("W8").Value = 3.7
$ws.Range("X8").Value = 2.14
$ws.Range("Y8").Value = 3.93
$ws.Range("Z8").Value = 2.02
$ws.Range("AA8").Value = 90.51000000000001
$ws.Range("AC8").Value = 1977
$ws.Range("AD8").Value = 10.39
$ws.Range("AE8").Value = 51097
$ws.Range("AF8").Value = 0.4
$ws.Range("AG8").Value = 600
$ws.Range("AH8").Value = 2.92
$ws.Range("AI8").Value = 30.35
$ws.Range("D9").Value = 26493
$ws.Range("E9").Value = 1087
$ws.Range("G9").Value = 875
$ws.Range("H9").Value = 644
$ws.Range("I9").Value = 659
$ws.Range("K9").Value = 27640
$ws.Range("L9").Value = 12804
$ws.Range("M9").Value = 14836
$ws.Range("N9").Value = 14788
$ws.Range("P9").Value = 1400
$ws.Range("Q9").Value = 1398
$ws.Range("R9").Value = -783
$ws.Range("S9").Value = -317
$ws.Range("T9").Value = 757
$ws.Range("U9").Value = 543
$ws.Range("W9").Value = 4.1
$ws.Range("X9").Value = 2.43
$ws.Range("Y9").Value = 4.53
$ws.Range("Z9").Value = 2.34
$ws.Range("AA9").Value = 86.31
$ws.Range("AC9").Value = 2351
$ws.Range("AD9").Value = 8.74
$ws.Range("AE9").Value = 52789
$ws.Range("AF9").Value = 0.39
$ws.Range("AG9").Value = 620
$ws.Range("AH9").Value = 3.02
$ws.Range("AI9").Value = 26.37
